$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad") from 45192 -> 45202 for existing rows 2..145
for ($r = 2; $r -le 145; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Excel recalculated/touched row 145's height explicitly when the row was
# last edited (matches surrounding rows which already carry ht="15").
$ws.Rows.Item(145).RowHeight = 15

# 2. Append the new avverkningsanmälan entry as row 146
$row = 146
$ws.Cells.Item($row, 1).Value = "A 46160-2023"        # A: Beteckning
$ws.Cells.Item($row, 2).Value = 45196                 # B: Datum
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 3).Value = 45202                 # C: Förändrad
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 4).Value = "UPPSALA LÄN"         # D: Län
$ws.Cells.Item($row, 5).Value = "KNIVSTA"             # E: Kommun
$ws.Cells.Item($row, 7).Value = 0.5                   # G: Area (ha)
$ws.Cells.Item($row, 8).Value = 0                     # H: Fridlysta
$ws.Cells.Item($row, 9).Value = 0                     # I: Signalarter
$ws.Cells.Item($row, 10).Value = 0                    # J: NT
$ws.Cells.Item($row, 11).Value = 0                    # K: VU
$ws.Cells.Item($row, 12).Value = 0                    # L: EN
$ws.Cells.Item($row, 13).Value = 0                    # M: CR
$ws.Cells.Item($row, 14).Value = 0                    # N: RE
$ws.Cells.Item($row, 15).Value = 0                    # O: Rödlistade
$ws.Cells.Item($row, 16).Value = 0                    # P: Hotade
$ws.Cells.Item($row, 17).Value = 0                    # Q: Alla arter
$ws.Cells.Item($row, 18).WrapText = $true              # R: Artnamn (styled, empty)
